# Insert a new, empty "Corps de texte" paragraph (left-indented one tab
# stop) right after the bookmark-terminated answer paragraph that ends in
# "...sortira courant novembre et donnera acces au site jusque a la fin de
# l'annee 2015." -- i.e. immediately after <w:bookmarkEnd w:id="0"/> and
# before the pre-existing blank paragraph that follows it.

$d = $word.ActiveDocument

# --- Locate the anchor paragraph robustly (no hard-coded paragraph index).
# "sortira" is unique to the paragraph that precedes the bookmark end.
$finder = $d.Content
$finder.Find.ClearFormatting()
$found = $finder.Find.Execute("sortira", $true, $false, $false, $false,
                               $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor text 'sortira' not found"
}
$anchorPos = $finder.Start

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($anchorPos -ge $candidate.Range.Start -and $anchorPos -lt $candidate.Range.End) {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not resolve anchor paragraph index"
}

# The paragraph right after the anchor is the pre-existing blank
# "Corps de texte" paragraph that immediately follows <w:bookmarkEnd/>.
# Collapsing its range to its own start and inserting a paragraph mark
# there creates the new paragraph in between, inheriting that following
# paragraph's formatting (style Corpsdetexte / spacing after=0 / jc=both).
$followingPara = $d.Paragraphs.Item($anchorIndex + 1)
$insertionPoint = $followingPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

# The freshly created paragraph now occupies the slot right after the
# anchor paragraph. Give it the one-tab-stop left indent (708 twips ==
# 35.4 points) that distinguishes it from the plain blank paragraph
# beneath it.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newPara.Format.LeftIndent = 35.4

# The automation layer leaves a formatting-only, textless run behind in
# the new paragraph. Typing a character and immediately deleting it
# flushes that stray run so the paragraph serializes with only <w:pPr>,
# matching a genuinely empty, manually-authored paragraph.
$cleanupRange = $newPara.Range
$cleanupRange.Collapse(1)
$cleanupRange.InsertBefore("x")
$d.Range($cleanupRange.Start, $cleanupRange.Start + 1).Delete()
